$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.746.65"
$ws.Range("E2").Value = "  +1.43%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.085.59"
$ws.Range("E3").Value = "  -0.17%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.15"
$ws.Range("E5").Value = "  -2.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.44"
$ws.Range("E6").Value = "  -0.78%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.079.73"
$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  +0.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.156"
$ws.Range("E10").Value = "  -1.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.24"
$ws.Range("E11").Value = "  -5.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  +1.24%  "

$ws.Range("E13").Value = "  +4.61%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.78"
$ws.Range("E14").Value = "  -0.55%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.590.72"
$ws.Range("E15").Value = "  -0.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.804.16"
$ws.Range("E16").Value = "  +1.30%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.091.04"
$ws.Range("E18").Value = "  -0.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.71"
$ws.Range("E19").Value = "  +0.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "488.78"
$ws.Range("E20").Value = "  -3.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.49"
$ws.Range("E21").Value = "  -0.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.703"
$ws.Range("E22").Value = "  -0.37%  "

$ws.Range("E23").Value = "  -0.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.74"
$ws.Range("E24").Value = "  +2.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.24"
$ws.Range("E25").Value = "  -0.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.73"
$ws.Range("E27").Value = "  -0.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.28"
$ws.Range("E28").Value = "  +0.68%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.27"
$ws.Range("E30").Value = "  +0.20%  "

$ws.Range("E31").Value = "  -2.07%  "

$ws.Range("E32").Value = "  +0.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.42"
$ws.Range("E33").Value = "  -3.86%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "57.45"
$ws.Range("E34").Value = "  -2.85%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.45"
$ws.Range("E35").Value = "  +5.06%  "

$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "495.21"
$ws.Range("E36").Value = "  -6.20%  "

$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.06"
$ws.Range("E37").Value = "  +2.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.265.84"
$ws.Range("E38").Value = "  +6.85%  "

$ws.Range("E39").Value = "  -2.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0802"
$ws.Range("E40").Value = "  +1.40%  "

$ws.Range("E41").Value = "  -1.26%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.16"
$ws.Range("E42").Value = "  +1.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.67"
$ws.Range("E43").Value = "  -0.62%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.258"
$ws.Range("E44").Value = "  +1.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.07"
$ws.Range("E46").Value = "  +1.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0537"
$ws.Range("E47").Value = "  +6.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "121.62"
$ws.Range("E48").Value = "  -0.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.85"
$ws.Range("E49").Value = "  +3.60%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.109"
$ws.Range("E50").Value = "  +2.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.36"
$ws.Range("E51").Value = "  +0.81%  "
